$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 196, pushing existing rows 196+ down to 197+.
# Copy row 196 (which currently holds the data that should become row 197)
# so the new row inherits the same formatting/values, then overwrite the
# three cells that differ (Fecha, Volumen, Origen).
$ws.Rows.Item(196).Copy() | Out-Null
$ws.Rows.Item(196).Insert() | Out-Null

$ws.Cells.Item(196, 4).Value = 44582
$ws.Cells.Item(196, 10).Value = 500
$ws.Cells.Item(196, 15).Value = "Provincia de Cautín"
